# Automation HUB process: insert additional yearly-report rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows before the current row 5, pushing the existing
# rows 5-7 down to become rows 9-11.
$ws.Range("A5:G8").EntireRow.Insert()

# New data for rows 5-8. Columns A, C, D, E, G look numeric/date-like and
# must stay stored as text (shared strings), so we force a text number
# format before assigning the value and then clear the format again so
# the cell keeps the default style.
$newRows = @(
    @{ Row = 5;  A = "646429"; B = "Waste management services"; C = "277397"; D = "55479.4";  E = "332876"; F = "CAD"; G = "2017-07-06" },
    @{ Row = 6;  A = "841881"; B = "Beverages and Catering";     C = "133356"; D = "26671.2";  E = "160027"; F = "CAD"; G = "2017-07-20" },
    @{ Row = 7;  A = "968494"; B = "Various paper supplies";     C = "170870"; D = "34174";    E = "205044"; F = "RON"; G = "2017-07-08" },
    @{ Row = 8;  A = "416313"; B = "Concierge Services";         C = "182562"; D = "36512.4";  E = "219074"; F = "RON"; G = "2017-08-24" }
)

foreach ($r in $newRows) {
    $rowRange = $ws.Range("A" + $r.Row + ":G" + $r.Row)
    $rowRange.NumberFormat = "@"

    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G

    $rowRange.ClearFormats()
}
